$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the very top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

$ws.Cells.Item(1, 1).Value = "email"
$ws.Cells.Item(1, 2).Value = "name"

$ws.Range("A2").Select()
